$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append the new bitcoin-buy row recorded on 2025-10-05.
# The leading apostrophe forces Excel to store the date-like string
# literally as text (matching the existing rows' inlineStr date cells)
# instead of auto-converting it to a date serial number; resetting the
# style back to Normal afterwards clears the transient quote-prefix
# formatting so the cell ends up with no explicit style, just like the
# other text-date cells in this column.
$ws.Range("A51").Value = "'10/05/2025"
$ws.Range("A51").Style = "Normal"

$ws.Range("B51").Value = 0.0003973100000000014
$ws.Range("C51").Value = 125846.3164783162
$ws.Range("D51").Value = 50
